$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '63.370.38'
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -4.56%  '
$c.Style = "Normal"

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.091.78'
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -4.76%  '
$c.Style = "Normal"

# Row 4
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  +0.01%  '
$c.Style = "Normal"

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '548.67'
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  -4.94%  '
$c.Style = "Normal"

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '137.27'
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -10.18%  '
$c.Style = "Normal"

# Row 7
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  +0.10%  '
$c.Style = "Normal"

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '3.083.81'
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  -4.78%  '
$c.Style = "Normal"

# Row 9
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -3.20%  '
$c.Style = "Normal"

# Row 10
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -4.66%  '
$c.Style = "Normal"

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '6.27'
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  -11.67%  '
$c.Style = "Normal"

# Row 12
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -3.40%  '
$c.Style = "Normal"

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '35.51'
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -5.81%  '
$c.Style = "Normal"

# Row 14
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -7.20%  '
$c.Style = "Normal"

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '3.588.31'
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  -4.72%  '
$c.Style = "Normal"

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '63.323.84'
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -4.58%  '
$c.Style = "Normal"

# Row 17
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -2.90%  '
$c.Style = "Normal"

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.088.30'
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -4.84%  '
$c.Style = "Normal"

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.75'
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -4.71%  '
$c.Style = "Normal"

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '489.37'
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -11.93%  '
$c.Style = "Normal"

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '13.64'
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -5.13%  '
$c.Style = "Normal"

# Row 22
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -3.02%  '
$c.Style = "Normal"

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '7.28'
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  -6.12%  '
$c.Style = "Normal"

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '79.08'
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -3.25%  '
$c.Style = "Normal"

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '12.38'
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -8.78%  '
$c.Style = "Normal"

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -0.29%  '
$c.Style = "Normal"

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '8.52'
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -8.08%  '
$c.Style = "Normal"

# Row 28
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  -6.04%  '
$c.Style = "Normal"

# Row 29
$ws.Range("B29").Value = 'FirstDigitalUSD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  +0.09%  '
$c.Style = "Normal"

# Row 30
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.97'
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -11.36%  '
$c.Style = "Normal"

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '26.66'
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -3.96%  '
$c.Style = "Normal"

# Row 32
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -4.20%  '
$c.Style = "Normal"

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '2.50'
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  -8.40%  '
$c.Style = "Normal"

# Row 34
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  +5.93%  '
$c.Style = "Normal"

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '508.88'
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -9.15%  '
$c.Style = "Normal"

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '6.05'
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  -4.93%  '
$c.Style = "Normal"

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '5.13'
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -10.05%  '
$c.Style = "Normal"

# Row 38
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -12.11%  '
$c.Style = "Normal"

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '3.152.23'
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -0.23%  '
$c.Style = "Normal"

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.0803'
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -6.67%  '
$c.Style = "Normal"

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.119'
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  -9.01%  '
$c.Style = "Normal"

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '8.17'
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -4.75%  '
$c.Style = "Normal"

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '2.64'
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -13.08%  '
$c.Style = "Normal"

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.258'
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -5.84%  '
$c.Style = "Normal"

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.07'
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  -9.93%  '
$c.Style = "Normal"

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '25.37'
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  -3.92%  '
$c.Style = "Normal"

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '120.82'
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  -3.97%  '
$c.Style = "Normal"

# Row 49
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -2.81%  '
$c.Style = "Normal"

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0₃0508'
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -8.26%  '
$c.Style = "Normal"

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '2.31'
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  +31.33%  '
$c.Style = "Normal"
